$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 96, shifting existing rows 96:181 down to 97:182
$ws.Rows.Item(96).Insert()

# Populate the newly inserted row 96 with the new data record
$ws.Range("A96").Value = 11
$ws.Range("B96").Value = "Vega Monumental Concepción"
$ws.Range("C96").Value = "Bíobío"
$ws.Range("D96").Value = 44729
$ws.Range("E96").Value = 8
$ws.Range("F96").Value = 100112003
$ws.Range("G96").Value = "Ajo"
$ws.Range("H96").Value = "Chino"
$ws.Range("I96").Value = "1a (guarda)"
$ws.Range("J96").Value = 310
$ws.Range("K96").Value = 17000
$ws.Range("L96").Value = 18000
$ws.Range("M96").Value = 17516
$ws.Range("N96").Value = "`$/caja 10 kilos"
$ws.Range("O96").Value = "China"
$ws.Range("P96").Value = 1752
$ws.Range("Q96").Value = 10
$ws.Range("R96").Value = "Hortaliza"
